# The workbook was opened in a newer Excel build, a value was entered into
# cell A1, and the worksheet's print orientation was set to Portrait.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the value 10 into A1
$ws.Range("A1").Value = 10

# Set the page orientation to Portrait
$ws.PageSetup.Orientation = $xlPortrait
